# Nalco aluminium ingot price sheet: prepend a new day's row (09-12-2025),
# shifting every existing data row down by one and carrying the last row's
# data into a brand new row at the bottom (126).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp
$newLastRow = $lastRow + 1

# Row $newLastRow (126) does not exist yet; give it the same cell formatting
# (style indexes) as the row currently being pushed off the bottom (125)
# before any values land in it, mirroring what Excel does when you drag a
# fill handle / copy a row down.
$ws.Range("A$lastRow`:F$lastRow").Copy() | Out-Null
$ws.Range("A$newLastRow`:F$newLastRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Columns that hold day-first (dd-mm-yyyy) date-looking text. Excel's normal
# text-entry parsing would silently reinterpret these as date serials, so we
# temporarily force Text format on the full working range before poking the
# values in, then restore General afterwards.
$dateColARange = $ws.Range("A2:A$newLastRow")
$dateColERange = $ws.Range("E2:E$newLastRow")
$dateColARange.NumberFormat = "@"
$dateColERange.NumberFormat = "@"

# Shift rows [2 .. lastRow] down to [3 .. lastRow+1], working bottom-up so a
# row is never overwritten before it has been read.
for ($r = $lastRow; $r -ge 2; $r--) {
    $dest = $r + 1
    $ws.Cells.Item($dest, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dest, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dest, 3).Value2 = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dest, 4).Value2 = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($dest, 5).Value2 = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($dest, 6).Value2 = $ws.Cells.Item($r, 6).Value2
}

# New top row: same price/circular info as the old row 2, new date.
$ws.Cells.Item(2, 1).Value2 = "09-12-2025"

$dateColARange.NumberFormat = "General"
$dateColERange.NumberFormat = "General"

# The hyperlink that used to live on F(lastRow) now also needs to cover the
# newly created F(newLastRow) cell (row 126), pointing at the same PDF.
$lastUrl = $ws.Cells.Item($newLastRow, 6).Value2
$ws.Hyperlinks.Add($ws.Cells.Item($newLastRow, 6), $lastUrl) | Out-Null

# Adding the hyperlink re-styles the cell with the default blue/underline
# "Hyperlink" look; put back the plain formatting used by every other link
# cell in this column so row 126 matches its neighbours.
$ws.Range("F$lastRow").Copy() | Out-Null
$ws.Range("F$newLastRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
